$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Login")

# New test case rows appended below existing data (rows 51-56), columns A:D
$newRows = @(
    @("Deals_Chat_ShipperUser_TC001",    "rogerdeals21+stan@gmail.com", "arewethere?", "Login successful"),
    @("Deals_Chat_ShipperUser_TC001(2)", "rogerdeals21+john@gmail.com", "arewethere?", "Login successful"),
    @("Deals_Chat_ShipperAdmin_TC002",   "rogerdeals21+nick@gmail.com", "arewethere?", "Login successful"),
    @("Deals_Chat_ShipperAdmin_TC002(2)","rogerdeals21+john@gmail.com", "arewethere?", "Login successful"),
    @("Deals_Chat_CarrierUser_TC003",    "rogerdeals21+john@gmail.com", "arewethere?", "Login successful"),
    @("Deals_Chat_CarrierUser_TC003(2)", "rogerdeals21+stan@gmail.com", "arewethere?", "Login successful")
)

$startRow = 51

# Shared strings end up ordered the way a human filled the sheet in the UI:
# column A for every new row first, then column B, then C, then D.
for ($col = 1; $col -le 4; $col++) {
    for ($i = 0; $i -lt $newRows.Count; $i++) {
        $r = $startRow + $i
        $ws.Cells.Item($r, $col).Value = $newRows[$i][$col - 1]
    }
}

# The first couple of new "Automation Test ID" cells picked up the
# vertically-centred formatting used throughout column A.
$ws.Cells.Item(51, 1).VerticalAlignment = -4108  # xlCenter
$ws.Cells.Item(52, 1).VerticalAlignment = -4108  # xlCenter

# Scroll/selection update to match the edited view
$ws.Application.Goto($ws.Range("A40"))
$ws.Range("A56").Select()
